$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-05-22 Thursday" "2025-05-23 Friday"

Replace-Text "327÷5=" "195÷9="
Replace-Text "886÷4=" "453÷8="
Replace-Text "539÷3=" "893÷4="
Replace-Text "271÷9=" "913÷2="
Replace-Text "283÷5=" "952÷2="

Replace-Text "920÷4=" "879÷5="
Replace-Text "745÷8=" "355÷3="
Replace-Text "943÷9=" "644÷9="
Replace-Text "369÷6=" "406÷3="
Replace-Text "733÷8=" "368÷5="

Replace-Text "623÷3=" "172÷8="
Replace-Text "628÷8=" "945÷4="
Replace-Text "745÷7=" "401÷8="
Replace-Text "662÷6=" "119÷4="
Replace-Text "447÷8=" "556÷8="

Replace-Text "910÷2=" "294÷9="
Replace-Text "817÷8=" "899÷8="
Replace-Text "492÷9=" "707÷4="
Replace-Text "627÷9=" "143÷6="
Replace-Text "649÷5=" "931÷5="

Replace-Text "399÷8=" "726÷2="
Replace-Text "803÷2=" "203÷2="
Replace-Text "351÷8=" "792÷2="
Replace-Text "665÷2=" "335÷9="
Replace-Text "611÷9=" "667÷3="
